# Added haptic feedback vibration motor control (NPN transistor on D10, driven from VBAT)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 44 for the NPN control pin details; everything below (the
# "Temperature" block etc.) shifts down by one row.
$ws.Rows.Item(44).Insert()

# Extend the blank separator row above the Vibration Motor section (row 40) so columns
# A-C also carry the plain Aptos Display spacer style, matching D40:O40.
$ws.Range("A40:C40").Font.Name = "Aptos Display"
$ws.Range("A40:C40").HorizontalAlignment = 1
$ws.Range("A40:C40").VerticalAlignment = -4107

# Fill the new row with the NPN transistor control pin entry (typed in D10, then the
# peripheral name, matching the order the new strings end up in the workbook).
$ws.Cells.Item(44, 2).Value = "D10"
$ws.Cells.Item(44, 1).Value = "Control (NPN)"

# The vibration motor is now switched through the NPN transistor, so it is wired to the
# battery rail (VBAT) instead of being fed straight off 3.3V.
$ws.Cells.Item(42, 2).Value = "VBAT"

$ws.Cells.Item(44, 3).Value = "connect to npn transistor"

# Re-style the GND row of the Vibration Motor block (row 43) so it matches the
# formatting of the "Vin" row directly above it (row 42), instead of its previous
# one-off styling.
$ws.Cells.Item(43, 1).Font.Name = "Aptos Display"
$ws.Cells.Item(43, 1).HorizontalAlignment = -4108
$ws.Cells.Item(43, 1).VerticalAlignment = -4107

$ws.Cells.Item(43, 3).Font.Name = "Aptos Display"
$ws.Cells.Item(43, 3).HorizontalAlignment = -4131
$ws.Cells.Item(43, 3).VerticalAlignment = -4107

# Style the new row: A44/B44 centered like the rest of the pin-name / pin-number
# columns, C44 left as the plain default-aligned note style.
$ws.Cells.Item(44, 1).Font.Name = "Aptos Display"
$ws.Cells.Item(44, 1).HorizontalAlignment = -4108
$ws.Cells.Item(44, 1).VerticalAlignment = -4108

$ws.Cells.Item(44, 2).Font.Name = "Aptos Display"
$ws.Cells.Item(44, 2).HorizontalAlignment = -4108
$ws.Cells.Item(44, 2).VerticalAlignment = -4108

$ws.Cells.Item(44, 3).Font.Name = "Aptos Display"
$ws.Cells.Item(44, 3).HorizontalAlignment = 1
$ws.Cells.Item(44, 3).VerticalAlignment = -4107

# Materialize the blank separator row (45) below the new entry, matching the plain
# Aptos Display font used for other un-styled spacer cells in the sheet.
$ws.Range("A45:C45").Font.Name = "Aptos Display"
$ws.Range("A45:C45").HorizontalAlignment = 1
$ws.Range("A45:C45").VerticalAlignment = -4107

# Restore the view/selection to where the author left off after adding the row.
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("E44").Select()
